$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.736.21'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '1.848.50'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("E6").Value = '  -0.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4647'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3856'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.78%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.77'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07906'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9700'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("D13").Value = '1.835.34'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.887'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.141'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("E16").Value = '  -0.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06610'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001028'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.44%  '

$ws.Range("D22").Value = '27.694.27'
$ws.Range("E22").Value = '  +0.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.358'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.285'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.23%  '

$ws.Range("D27").Value = '2.045.81'
$ws.Range("E27").Value = '  -1.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.071'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.340'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.58%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09419'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9446'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.586'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.262'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.329'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05995'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02202'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.230'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.004'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.154'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5796'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1839'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.278'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5441'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.88'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.925'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06843'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -32.60%  '
